$d = $word.ActiveDocument

# Start from the range of the last (currently only) paragraph: "Hello there this is my new doc file"
$lastRange = $d.Paragraphs.Last.Range

# Insert a new paragraph mark right after it, then fill the new paragraph with text
$lastRange.InsertParagraphAfter()

# The newly created paragraph is now the last paragraph in the document
$newPara = $d.Paragraphs.Last
$newPara.Range.InsertAfter("This is changed")

# Match the formatting of the preceding paragraph: 36pt (w:sz 72 half-points) font, en-IN language
$newPara.Range.Font.Size = 36
$newPara.Range.LanguageID = "en-IN"
